$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row (4) describing the "Raw Strings in R" code snippet.
# Values are entered in the order title, url, tags, src, description so
# that new shared-string entries are appended in that same order.
$ws.Range("A4").Value = "Raw Strings in R"
$ws.Range("C4").Value = "https://sciencificity.github.io/raw-strings-r/"
$ws.Range("E4").Value = "R; Raw Strings in R"
$ws.Range("B4").Value = "images/mae-mu--dyxcGiP-rE-unsplash1.jpg"
$ws.Range("D4").Value = 'In Version 4.0.0 of R raw strings were added r"(...)"'
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

# Add a hyperlink on the url cell, mirroring the existing one on C2, and
# restore the shared "Hyperlink" cell style afterwards.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://sciencificity.github.io/raw-strings-r/") | Out-Null
$ws.Range("C4").Style = "Hyperlink"

# Leave the active selection on the new description cell, matching the
# state the workbook was saved in.
$ws.Range("D4").Select()
